# B1--and-B2-PowerPoint.pptx edit
#
# 1) Table on slide 5 switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 1" table style.
# 2) The deck's theme colour scheme (the slide master's theme, stored as
#    ppt/theme/theme2.xml) is swapped from the "Integral / Red Violet"
#    palette back to the plain "Office Theme" palette (the palette that
#    used to live in ppt/theme/theme1.xml, the Notes Master's theme).
#    Font scheme / format scheme are identical between the two themes,
#    so only the 12 theme colours need to change.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ------------------------------------------------
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B532372D-DC2D-42B6-99F3-D98F25C9C58D}")
    }
}

# --- 2. Theme colour scheme swap (Integral -> Office Theme) ------------------
# RGB() packs as 0x00BBGGRR, matching classic VBA colour literals.
$officeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A
    4  = 0xE6E6E7   # lt2      E7E6E6
    5  = 0xD59B5B   # accent1  5B9BD5
    6  = 0x317DED   # accent2  ED7D31
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000
    9  = 0xC47244   # accent5  4472C4
    10 = 0x47AD70   # accent6  70AD47
    11 = 0xC16305   # hlink    0563C1
    12 = 0x724F95   # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i]
}
